$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: I1 = "I0", J1 = "IF" (copy formatting/style from H1, the last existing header) ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2..6: I = 1 (constant), J = same numeric value as column H ---
for ($r = 2; $r -le 6; $r++) {
    $hVal = $ws.Range("H$r").Value2
    $ws.Range("I$r").Value = 1
    $ws.Range("J$r").Value = $hVal
}
